# Fix pq for Chinese and other small bugs
# - Insert a new "level of education" question as row 4.
# - Clarify the "Do you read in any other language(s) than ..." question text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a brand new row above the current row 4 ("Years of education" question
#    stays at row 3; everything from old row 4 onward shifts down by one).
$ws.Rows.Item(4).Insert()

# 2) Fill in the new question row (row 4): "2a - Please report the highest level
#    of education you have completed" with its answer options.
$ws.Range("A4").Value = "2a"
$ws.Range("B4").Value = "Please report the highest level of education you have completed"
$ws.Range("C4").Value = "Less than primary education"
$ws.Range("D4").Value = "Primary education"
$ws.Range("E4").Value = "Lower secondary education"
$ws.Range("F4").Value = "Upper secondary education"
$ws.Range("G4").Value = "Post-secondary non-tertiary education"
$ws.Range("H4").Value = "Short-cycle tertiary education"
$ws.Range("I4").Value = "Bachelor&apos;s or equivalent"
$ws.Range("J4").Value = "Master&apos;s or equivalent"
$ws.Range("K4").Value = "Doctoral or equivalent"
$ws.Range("M4").Value = 1
$ws.Range("N4").Value = "level_education"
$ws.Range("O4").Value = "dropdown"

# Formatting used on this new row: most cells in Arial/black, L4 & O4 in
# Aptos Narrow/black (matches the rest of the "answer type" column style).
$ws.Range("A4:K4").Font.Name = "Arial"
$ws.Range("A4:K4").Font.Color = 0
$ws.Range("M4:N4").Font.Name = "Arial"
$ws.Range("M4:N4").Font.Color = 0
$ws.Range("L4").Font.Color = 0
$ws.Range("O4").Font.Color = 0

# 3) The old row 21 ("Do you read in any other language(s) than ...") is now row 22
#    after the insertion above. Make the question explicit about which language.
$ws.Range("B22").Value = "Do you read in any other language(s) than [insert current language]"
$ws.Range("B22").Font.Color = 0
